# add preview in datasets
# Insert two new metadata rows (dataset.preview.table / dataset.preview.line)
# right above "dataset.commit.id", each holding a multi-line OJS-style
# preview query used by the dataset viewer.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("metadata")

# Make room for the two new rows (old row 4 -> new row 6, etc.)
$ws.Rows("4:5").Insert()

$tableFormula = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nlimit(start:0, length:5);"
$lineFormula  = "source(ds:'{{dataset.id}}');`nquery([`n  { dim:'time', role:'row', items:[] },`n  { dim:'indicator', role:'col', items:[] } `n]);`nformat(p:3);`norder(dir:'row', index:-1, asc:'az');`nline(x:-1);"

$ws.Range("A4").Value = "dataset.preview.table"
$ws.Range("B4").Value = $tableFormula

$ws.Range("A5").Value = "dataset.preview.line"
$ws.Range("B5").Value = $lineFormula

# Match the wrapped / vertically-centered style already used for other
# multi-line cells, and size the rows to show the whole snippet.
$ws.Range("A4:B5").WrapText = $true
$ws.Rows("4:5").RowHeight = 120

# Reset the selection: B9 is the row that used to be B7 before the insert.
$ws.Range("B9").Select()
